$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(6)
$r = $p.Range
$found0 = $r.Find.Execute("An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io", $true, $false, $false, $false, $false, $true, 1, $false, "Một email được gửi đến các đối tác ở quốc gia mục tiêu đã phản hồi đồng ý nhưng không nộp tài liệu của họ trước hạn chót. Chúng tôi sẽ thu hồi lời mời của họ. It will be sent via customer.io", 2)
Write-Output "op0 (para 6) found=$found0"

$p = $d.Paragraphs.Item(9)
$r = $p.Range
$found1 = $r.Find.Execute("Invited partners who didn’t submit their documents on time", $true, $false, $false, $false, $false, $true, 1, $false, "Các đối tác được mời chưa nộp tài liệu của họ đúng hạn", 2)
Write-Output "op1 (para 9) found=$found1"

$p = $d.Paragraphs.Item(12)
$r = $p.Range
$found2 = $r.Find.Execute(": Your ", $true, $false, $false, $false, $false, $true, 1, $false, ": Đăng ký sự kiện ", 2)
Write-Output "op2 (para 12) found=$found2"

$p = $d.Paragraphs.Item(12)
$r = $p.Range
$found3 = $r.Find.Execute(" registration", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Output "op3 (para 12) found=$found3"

$p = $d.Paragraphs.Item(14)
$r = $p.Range
$found4 = $r.Find.Execute("We didn’t receive your documents on time", $true, $false, $false, $false, $false, $true, 1, $false, "Chúng tôi không nhận được giấy tờ của bạn đúng hạn", 2)
Write-Output "op4 (para 14) found=$found4"

$p = $d.Paragraphs.Item(16)
$r = $p.Range
$found5 = $r.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "Xin chào ", 2)
Write-Output "op5 (para 16) found=$found5"

$p = $d.Paragraphs.Item(17)
$r = $p.Range
$found6 = $r.Find.Execute("We didn’t receive your documents by the deadline (", $true, $false, $false, $false, $false, $true, 1, $false, "Chúng tôi không nhận được giấy tờ bạn cần cung cấp cho chúng tôi trước thời hạn (", 2)
Write-Output "op6 (para 17) found=$found6"

$p = $d.Paragraphs.Item(17)
$r = $p.Range
$found7 = $r.Find.Execute("). Unfortunately, we’re unable to proceed with your registration for the ", $true, $false, $false, $false, $false, $true, 1, $false, "). Vì vậy rất tiếc, chúng tôi không thể tiếp tục xử lý đơn đăng ký của bạn cho sự kiện ", 2)
Write-Output "op7 (para 17) found=$found7"

$p = $d.Paragraphs.Item(19)
$r = $p.Range
$found8 = $r.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua ", 2)
Write-Output "op8 (para 19) found=$found8"

# " or " -> " hoặc " (narrowed match to avoid hyperlink rPr bleed; run sits between two hyperlinks)
$p = $d.Paragraphs.Item(19)
$r = $p.Range
$found9 = $r.Find.Execute("or", $true, $true, $false, $false, $false, $true, 1, $false, "hoặc", 2)
Write-Output "op9 (para 19) found=$found9"

$p = $d.Paragraphs.Item(20)
$r = $p.Range
$found10 = $r.Find.Execute("If you have any questions, please contact your country manager, ", $true, $false, $false, $false, $false, $true, 1, $false, "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn ", 2)
Write-Output "op10 (para 20) found=$found10"

$p = $d.Paragraphs.Item(20)
$r = $p.Range
$found11 = $r.Find.Execute(", at ", $true, $false, $false, $false, $false, $true, 1, $false, ", qua email ", 2)
Write-Output "op11 (para 20) found=$found11"

$p = $d.Paragraphs.Item(20)
$r = $p.Range
$found12 = $r.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " hoặc số ", 2)
Write-Output "op12 (para 20) found=$found12"
